$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 52 - Panama / Kuna Yala / Aligandi
$ws.Range("A52").Value = "Panama"
$ws.Range("B52").Value = 9.134
$ws.Range("C52").Value = -77.951
$ws.Range("D52").Value = "-"
$ws.Range("E52").Value = "1842 (2008)"
$ws.Range("F52").Value = "Aligandi"
$ws.Range("G52").Value = "Kuna Yala"
$ws.Range("H52").Value = 42075
$ws.Range("H52").NumberFormat = "yyyy/mm/dd;@"
$ws.Range("I52").Value = "Zika"
$ws.Range("J52").Value = "First 3 cases in the country (20-40 years in age)"
$ws.Range("K52").Value = "http://www.prensa.com/sociedad/Detectan-casos-virus-zika-Panama_0_4360814004.html"

# Row 53 - Guatemala / Zacapa
$ws.Range("A53").Value = "Guatemala"
$ws.Range("B53").Value = 15.036
$ws.Range("C53").Value = -89.768
$ws.Range("D53").Value = "-"
$ws.Range("E53").Value = "-"
$ws.Range("F53").Value = "-"
$ws.Range("G53").Value = "Zacapa"
$ws.Range("H53").Value = "2015-24-11"
$ws.Range("I53").Value = "Zika"
$ws.Range("J53").Value = "22 cases"
$ws.Range("K53").Value = "http://informaciontotal.com.mx/2015-12-09-ea732713/guatemala-confirma-29-casos-de-virus-del-zika/"

# Row 54 - Guatemala / Escuintla
$ws.Range("A54").Value = "Guatemala"
$ws.Range("B54").Value = 14.194
$ws.Range("C54").Value = -91.294
$ws.Range("D54").Value = "-"
$ws.Range("E54").Value = "-"
$ws.Range("F54").Value = "-"
$ws.Range("G54").Value = "Escuintla"
$ws.Range("H54").Value = "2015-24-11"
$ws.Range("I54").Value = "Zika"
$ws.Range("J54").Value = "7 cases"
$ws.Range("K54").Value = "http://informaciontotal.com.mx/2015-12-09-ea732713/guatemala-confirma-29-casos-de-virus-del-zika/"

# Row 55 - Panama (misspelled "Pnanama") / Guna Yala
$ws.Range("A55").Value = "Pnanama"
$ws.Range("B55").Value = 9.06
$ws.Range("C55").Value = -78.881
$ws.Range("D55").Value = "-"
$ws.Range("E55").Value = "15 541 (2000)"
$ws.Range("F55").Value = "-"
$ws.Range("G55").Value = "Guna Yala"
$ws.Range("H55").Value = 42320
$ws.Range("H55").NumberFormat = "yyyy/mm/dd;@"
$ws.Range("I55").Value = "Zika"
$ws.Range("J55").Value = "1 new case in country bringing total to 4"
$ws.Range("K55").Value = "http://laestrella.com.pa/panama/nacional/minsa-confirma-nuevo-caso-virus-zika-panama/23909238"

$ws.Range("C56").Select()
